$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the typo "Etobiceto" -> "Etobicoke" in the cluster-name column (B),
# e.g. "North Etobicoke & South Etobiceto" -> "North Etobicoke & South Etobicoke"
$old = "North Etobicoke & South Etobiceto"
$new = "North Etobicoke & South Etobicoke"

$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq $old) {
        $cell.Value2 = $new
    }
}

# Move the selection/scroll position to reflect where the edit was made (cosmetic).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 149
$ws.Range("B149").Select()
